# Daily attendance processing - 2025-10-13 08:53:55
#
# For every "Recorded By" cell in column G (row 2 through the last used
# row) that holds a comma-separated list of recorders, rotate the list so
# the last entry moves to the front (everything else keeps its relative
# order). Cells with a single recorder (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -ne $null -and $v -ne "" -and $v -like "*,*") {
        $parts = $v -split ", "
        $n = $parts.Length
        $last = $parts[$n - 1]
        $rest = $parts[0..($n - 2)]
        $rotated = @($last) + $rest
        $cell.Value2 = $rotated -join ", "
    }
}
